$d = $word.ActiveDocument

# Locate the paragraph that holds the Google Drive hyperlink ("Video link:-" is the
# paragraph right before it) so we can insert the two new "Git hub file link:" /
# URL paragraphs immediately after it, before the existing blank paragraph.
$driveRange = $d.Content
$found = $driveRange.Find.Execute(
    "https://drive.google.com/file/d/1Jt5wDwVVV6fXnkShqmJJaGXii2SIYiE5/view?usp=sharing",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$hyperlinkPara = $driveRange.Paragraphs.Last
$insertionPoint = $d.Range($hyperlinkPara.Range.End, $hyperlinkPara.Range.End)

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="40"/>
<w:szCs w:val="40"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="40"/>
<w:szCs w:val="40"/>
</w:rPr>
<w:t>Git hub file link:</w:t>
</w:r>
</w:p>
<w:p>
<w:pPr>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="40"/>
<w:szCs w:val="40"/>
<w:u w:val="single"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="40"/>
<w:szCs w:val="40"/>
<w:u w:val="single"/>
</w:rPr>
<w:t>https://github.com/Debasis-Behera4143/zigzag-matrix-task.git</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertionPoint.InsertXML($xml) | Out-Null
